# Adding 3-may data and modifications in model, to separate Madrid,
# La Rioja and Pais Vasco from others.
#
# Appends the "tabula-Actualizacion_95_COVID-19" extraction (2020-05-03,
# rows 851-870) to the bottom of the Data sheet, following the exact
# layout/format used by the previous day's block (rows 841-850), and
# registers the matching named range used by the Tabula PDF extraction
# add-in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$date = 43954   # 2020-05-03

# row, gender, age group, casos, hospit, pct_hosp, ingr_UCI, pct_uci, deaths, letalidad, pct_death
$rows = @(
    @(851, "Fem",    "0-9",    330,   107,  0.3,  13,    0.6,  1,    0,    0.3),
    @(852, "Fem",    "10-19",  729,   131,  0.3,  8,     0.3,  3,    0,    0.4),
    @(853, "Fem",    "20-29",  8069,  727,  1.9,  38,    1.6,  8,    0.1,  0.1),
    @(854, "Fem",    "30-39",  12791, 1726, 4.5,  94,    4.1,  22,   0.3,  0.2),
    @(855, "Fem",    "40-49",  18558, 3453, 9,    223,   9.6,  70,   1,    0.4),
    @(856, "Fem",    "50-59",  21932, 5529, 14.5, 432,   18.6, 161,  2.2,  0.7),
    @(857, "Fem",    "60-69",  15069, 6591, 17.3, 679,   29.3, 439,  6.1,  2.9),
    @(858, "Fem",    "70-79",  13537, 8098, 21.2, 664,   28.7, 1352, 18.7, 10),
    @(859, "Fem",    "80-89",  19462, 8375, 21.9, 140,   6,    3089, 42.7, 15.9),
    @(860, "Fem",    "90 y +", 11355, 3446, 9,    26,    1.1,  2093, 28.9, 18.4),
    @(861, "Masc",   "0-9",    401,   140,  0.3,  25,    0.5,  1,    0,    0.2),
    @(862, "Masc",   "10-19",  589,   123,  0.2,  11,    0.2,  2,    0,    0.3),
    @(863, "Masc",   "20-29",  3953,  657,  1.3,  50,    1,    14,   0.1,  0.4),
    @(864, "Masc",   "30-39",  7521,  1852, 3.8,  162,   3.1,  35,   0.4,  0.5),
    @(865, "Masc",   "40-49",  13010, 4859, 9.9,  495,   9.4,  113,  1.1,  0.9),
    @(866, "Masc",   "50-59",  16788, 8029, 16.3, 1062,  20.3, 405,  4.1,  2.4),
    @(867, "Masc",   "60-69",  16698, 10137,20.6, 1744,  33.3, 1079, 10.9, 6.5),
    @(868, "Masc",   "70-79",  17046, 12037,24.4, 1502,  28.6, 2908, 29.3, 17.1),
    @(869, "Masc",   "80-89",  14480, 9206, 18.7, 172,   3.3,  3997, 40.2, 27.6),
    @(870, "Masc",   "90 y +", 4527,  2236, 4.5,  21,    0.4,  1386, 13.9, 30.6)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $gender = $r[1]
    $age    = $r[2]
    $casos  = $r[3]
    $hospit = $r[4]
    $pctHosp= $r[5]
    $uci    = $r[6]
    $pctUci = $r[7]
    $deaths = $r[8]
    $letal  = $r[9]
    $pctDeath = $r[10]

    $ws.Range("A$rowNum").Value = $date
    $ws.Range("A$rowNum").NumberFormat = "d-mmm"

    $ws.Range("B$rowNum").Value = $gender

    $ws.Range("C$rowNum").Value = $age
    if ($age -eq "10-19") {
        $ws.Range("C$rowNum").NumberFormat = "@"
    }

    $ws.Range("D$rowNum").Value = $casos
    $ws.Range("E$rowNum").Value = $hospit
    $ws.Range("F$rowNum").Value = $pctHosp
    $ws.Range("G$rowNum").Value = $uci
    $ws.Range("H$rowNum").Value = $pctUci
    $ws.Range("I$rowNum").Value = $deaths
    $ws.Range("J$rowNum").Value = $letal
    $ws.Range("K$rowNum").Value = $pctDeath

    foreach ($col in @("D", "E", "G", "I")) {
        $cell = $ws.Range("$col$rowNum")
        if ($cell.Value -ge 1000) {
            $cell.NumberFormat = "#,##0"
        }
    }
}

# New defined name (localSheetId=0 / scoped to "Data") registered for the
# Tabula extraction, mirroring the pattern of tabula_Actualizacion_94_COVID_19.
$rng = $ws.Range("C851:K870")
$ws.Names.Add("tabula_Actualizacion_95_COVID_19", $rng)

# Reflect the new bottom of the data range / scrolled viewport, matching
# what Excel records after appending + scrolling to the new last row.
$ws.Activate()
$win = $excel.ActiveWindow
if ($win) {
    $win.ScrollRow = 835
}
$ws.Range("B870").Select()
